# Swap the full record (all columns except the row-index column A)
# between each of the following pairs of rows on the active sheet.
# This matches a source re-sort where two neighboring fixtures had
# their row order flipped while the running "id" index (column A)
# stayed attached to the row position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

$rowPairs = @(
    @(39, 40),
    @(74, 75),
    @(111, 112),
    @(186, 187),
    @(243, 244),
    @(247, 248),
    @(256, 257),
    @(260, 261),
    @(278, 279),
    @(293, 294),
    @(300, 301),
    @(303, 304)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    foreach ($col in $cols) {
        $cell1 = $ws.Range("$col$r1")
        $cell2 = $ws.Range("$col$r2")

        $v1 = $cell1.Value2
        $v2 = $cell2.Value2

        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}
